$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56
$ws.Range("B56").Value = 'schema_types_5_01'
$ws.Range("C56").Value = 'Passed'
$ws.Range("D56").Value = 'DFDL-5-015R'
$ws.Range("E56").Value = 'High'
$ws.Range("F56").Value = 3
$ws.Range("G56").Value = 'text'
$ws.Range("H56").Value = 'DFDL-219'

# Row 57
$ws.Range("B57").Value = 'schema_types_5_05'
$ws.Range("C57").Value = 'Passed'
$ws.Range("D57").Value = 'DFDL-5-020R'
$ws.Range("E57").Value = 'High'
$ws.Range("F57").Value = 3
$ws.Range("G57").Value = 'text'
$ws.Range("H57").Value = 'DFDL-219'

# Row 58
$ws.Range("B58").Value = 'syntax_entities_6_01'
$ws.Range("C58").Value = 'Passed'
$ws.Range("D58").Value = 'DFDL-6-041R'
$ws.Range("E58").Value = 'High'
$ws.Range("F58").Value = 3
$ws.Range("G58").Value = 'byte'
$ws.Range("H58").Value = 'DFDL-219'

# Row 59
$ws.Range("B59").Value = 'syntax_entities_6_03'
$ws.Range("C59").Value = 'Passed'
$ws.Range("D59").Value = 'DFDL-6-045R'
$ws.Range("E59").Value = 'High'
$ws.Range("F59").Value = 3
$ws.Range("G59").Value = 'byte'
$ws.Range("H59").Value = 'DFDL-219'

# Row 60
$ws.Range("B60").Value = 'NextLine'
$ws.Range("C60").Value = 'Passed'
$ws.Range("D60").Value = 'DFDL-6-045R'
$ws.Range("E60").Value = 'High'
$ws.Range("F60").Value = 3
$ws.Range("G60").Value = 'text'
$ws.Range("H60").Value = 'DFDL-219'

# Row 61
$ws.Range("B61").Value = 'LineSeparator'
$ws.Range("C61").Value = 'Passed'
$ws.Range("D61").Value = 'DFDL-6-045R'
$ws.Range("E61").Value = 'High'
$ws.Range("F61").Value = 3
$ws.Range("G61").Value = 'text'
$ws.Range("H61").Value = 'DFDL-219'

# Row 62
$ws.Range("B62").Value = 'LineFeed'
$ws.Range("C62").Value = 'Passed'
$ws.Range("D62").Value = 'DFDL-6-045R'
$ws.Range("E62").Value = 'High'
$ws.Range("F62").Value = 3
$ws.Range("G62").Value = 'text'
$ws.Range("H62").Value = 'DFDL-219'

# Row 63
$ws.Range("B63").Value = 'CarriageReturn'
$ws.Range("C63").Value = 'Passed'
$ws.Range("D63").Value = 'DFDL-6-045R'
$ws.Range("E63").Value = 'High'
$ws.Range("F63").Value = 3
$ws.Range("G63").Value = 'text'
$ws.Range("H63").Value = 'DFDL-219'

# Row 64
$ws.Range("B64").Value = 'FormFeed'
$ws.Range("C64").Value = 'Passed'
$ws.Range("D64").Value = 'DFDL-6-042R'
$ws.Range("E64").Value = 'High'
$ws.Range("F64").Value = 3
$ws.Range("G64").Value = 'byte'
$ws.Range("H64").Value = 'DFDL-219'

# Row 65
$ws.Range("B65").Value = 'CarriageReturn_byte'
$ws.Range("C65").Value = 'Passed'
$ws.Range("D65").Value = 'DFDL-6-045R'
$ws.Range("E65").Value = 'High'
$ws.Range("F65").Value = 3
$ws.Range("G65").Value = 'byte'
$ws.Range("H65").Value = 'DFDL-141'

# Row 66
$ws.Range("B66").Value = 'LineFeed_byte'
$ws.Range("C66").Value = 'Passed'
$ws.Range("D66").Value = 'DFDL-6-045R'
$ws.Range("E66").Value = 'High'
$ws.Range("F66").Value = 3
$ws.Range("G66").Value = 'byte'
$ws.Range("H66").Value = 'DFDL-141'

# Row 67
$ws.Range("B67").Value = 'LineSeparator_byte'
$ws.Range("C67").Value = 'Passed'
$ws.Range("D67").Value = 'DFDL-6-045R'
$ws.Range("E67").Value = 'High'
$ws.Range("F67").Value = 3
$ws.Range("G67").Value = 'byte'
$ws.Range("H67").Value = 'DFDL-141'

# Row 68
$ws.Range("B68").Value = 'NextLine_byte'
$ws.Range("C68").Value = 'Passed'
$ws.Range("D68").Value = 'DFDL-6-045R'
$ws.Range("E68").Value = 'High'
$ws.Range("F68").Value = 3
$ws.Range("G68").Value = 'byte'
$ws.Range("H68").Value = 'DFDL-141'

# Row 69
$ws.Range("B69").Value = 'CRLF_byte'
$ws.Range("C69").Value = 'Passed'
$ws.Range("D69").Value = 'DFDL-6-045R'
$ws.Range("E69").Value = 'High'
$ws.Range("F69").Value = 3
$ws.Range("G69").Value = 'byte'
$ws.Range("H69").Value = 'DFDL-141'

# Row 70
$ws.Range("B70").Value = 'lengthKindDelimited_01'
$ws.Range("C70").Value = 'Passed'
$ws.Range("D70").Value = 'DFDL-12-048R'
$ws.Range("E70").Value = 'High'
$ws.Range("F70").Value = 3
$ws.Range("G70").Value = 'text'
$ws.Range("H70").Value = 'DFDL-142'

# Row 71
$ws.Range("B71").Value = 'lengthKindDelimited_02'
$ws.Range("C71").Value = 'Passed'
$ws.Range("D71").Value = 'DFDL-12-048R'
$ws.Range("E71").Value = 'High'
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 'text'
$ws.Range("H71").Value = 'DFDL-142'

# Row 72
$ws.Range("B72").Value = 'delimiter_12_03'
$ws.Range("C72").Value = 'Passed'
$ws.Range("D72").Value = 'DFDL-12-032R'
$ws.Range("E72").Value = 'High'
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 'text'
$ws.Range("H72").Value = 'DFDL-109'

# Row 73
$ws.Range("B73").Value = 'SeqGrp_01'
$ws.Range("C73").Value = 'Passed'
$ws.Range("D73").Value = 'DFDL-14-008R'
$ws.Range("E73").Value = 'High'
$ws.Range("F73").Value = 3
$ws.Range("G73").Value = 'text'
$ws.Range("H73").Value = 'DFDL-109'

# Row 74
$ws.Range("B74").Value = 'litNil3'
$ws.Range("C74").Value = 'Passed'
$ws.Range("D74").Value = 'DFDL-13-236R'
$ws.Range("E74").Value = 'High'
$ws.Range("F74").Value = 3
$ws.Range("G74").Value = 'text'
$ws.Range("H74").Value = 'DFDL-199'

# Row 75
$ws.Range("B75").Value = 'DelimProp_01'
$ws.Range("C75").Value = 'Passed'
$ws.Range("D75").Value = 'DFDL-12-032R'
$ws.Range("E75").Value = 'High'
$ws.Range("F75").Value = 3
$ws.Range("G75").Value = 'text'
$ws.Range("H75").Value = 'DFDL-203'

# Row 76
$ws.Range("B76").Value = 'ParseSequence4'
$ws.Range("C76").Value = 'Passed'
$ws.Range("D76").Value = 'DFDL-12-032R'
$ws.Range("E76").Value = 'High'
$ws.Range("F76").Value = 3
$ws.Range("G76").Value = 'text'
$ws.Range("H76").Value = 'DFDL-203'

# Row 77
$ws.Range("B77").Value = 'ParseSequence5'
$ws.Range("C77").Value = 'Passed'
$ws.Range("D77").Value = 'DFDL-12-032R'
$ws.Range("E77").Value = 'High'
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 'text'
$ws.Range("H77").Value = 'DFDL-203'

# Update final selection to match target view state
$ws.Range("H76").Select()

